$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New correlation values for row 2 (and mirrored column B) of the matrix.
$ws.Range("C2").Value = -0.7047556261015829
$ws.Range("D2").Value = -0.04743437747629582
$ws.Range("E2").Value = 0.02275530254629318
$ws.Range("G2").Value = 0.01984950350403698
$ws.Range("H2").Value = -0.07735816559944704
$ws.Range("I2").Value = -0.1403061574535686
$ws.Range("J2").Value = -0.04671663250852448
$ws.Range("K2").Value = 0.02676678503401862
$ws.Range("L2").Value = -0.07828834604569371
$ws.Range("M2").Value = 0.004732017893733165
$ws.Range("N2").Value = -0.08118141911519945
$ws.Range("O2").Value = -0.02673676715644199

# Mirrored values in column B (matrix symmetry)
$ws.Range("B3").Value = -0.7047556261015829
$ws.Range("B4").Value = -0.04743437747629582
$ws.Range("B5").Value = 0.02275530254629318
$ws.Range("B7").Value = 0.01984950350403698
$ws.Range("B8").Value = -0.07735816559944704
$ws.Range("B9").Value = -0.1403061574535686
$ws.Range("B10").Value = -0.04671663250852448
$ws.Range("B11").Value = 0.02676678503401862
$ws.Range("B12").Value = -0.07828834604569371
$ws.Range("B13").Value = 0.004732017893733165
$ws.Range("B14").Value = -0.08118141911519945
$ws.Range("B15").Value = -0.02673676715644199
